$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new status note in column C next to the "Status" row (row 27)
$ws.Range("C27").Value = "Just programmed it."

# Move the active selection to C28, matching the post-edit cursor position
$ws.Range("C28").Select()
